$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5000
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5568

$ws.Range("H74").Value = 2699.923
$ws.Range("I74").Value = 2716.5833
$ws.Range("K74").Value = 2716.5833
$ws.Range("M74").Value = -1780.5833

$ws.Range("H77").Value = 2699.923
$ws.Range("I77").Value = 2716.5833
$ws.Range("K77").Value = 13582.9165
$ws.Range("M77").Value = -8902.916499999999

$ws.Range("H95").Value = 38333.332
$ws.Range("J95").Value = 38333.332
$ws.Range("L95").Value = 38333.332
$ws.Range("N95").Value = -43825.332

$ws.Range("H105").Value = 53333
$ws.Range("J105").Value = 53333
$ws.Range("L105").Value = 53333
$ws.Range("N105").Value = -60321

$ws.Range("H120").Value = 50001
$ws.Range("J120").Value = 50001
$ws.Range("L120").Value = 50001
$ws.Range("N120").Value = -59677

$ws.Range("H138").Value = 3795.8823
$ws.Range("I138").Value = 2623.6191
$ws.Range("J138").Value = 5689.5386
$ws.Range("K138").Value = 7870.8573
$ws.Range("L138").Value = 17068.6158
$ws.Range("M138").Value = -2730.8573
$ws.Range("N138").Value = -27348.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6628.0586
$ws.Range("I63").Value = 1838.5
$ws.Range("J63").Value = 7266.6665
$ws.Range("K63").Value = 1838.5
$ws.Range("L63").Value = 7266.6665
$ws.Range("M63").Value = -1152.5
$ws.Range("N63").Value = -8638.666499999999

$ws.Range("H66").Value = 6628.0586
$ws.Range("I66").Value = 1838.5
$ws.Range("J66").Value = 7266.6665
$ws.Range("K66").Value = 9192.5
$ws.Range("L66").Value = 36333.3325
$ws.Range("M66").Value = -5760.5
$ws.Range("N66").Value = -43197.3325

$ws.Range("H97").Value = 656.6667
$ws.Range("I97").Value = 656.6667
$ws.Range("K97").Value = 656.6667
$ws.Range("M97").Value = -160.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 6250441
$ws.Range("I7").Value = 7143218.5
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 7143218.5
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -7143105.5
$ws.Range("N7").Value = -1226

$ws.Range("H20").Value = 13249.5
$ws.Range("I20").Value = 11499
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 11499
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -11252
$ws.Range("N20").Value = -15494

$ws.Range("H36").Value = 2000
$ws.Range("I36").Value = 2000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1466

$ws.Range("H64").Value = 1290
$ws.Range("I64").Value = 1187.5
$ws.Range("J64").Value = 1495
$ws.Range("K64").Value = 1187.5
$ws.Range("L64").Value = 1495
$ws.Range("M64").Value = -962.5
$ws.Range("N64").Value = -1945

$ws.Range("H67").Value = 1290
$ws.Range("I67").Value = 1187.5
$ws.Range("J67").Value = 1495
$ws.Range("K67").Value = 1187.5
$ws.Range("L67").Value = 1495
$ws.Range("M67").Value = -407.5
$ws.Range("N67").Value = -3055

$ws.Range("H99").Value = 1816.875
$ws.Range("I99").Value = 1833.5714
$ws.Range("J99").Value = 1700
$ws.Range("K99").Value = 1833.5714
$ws.Range("L99").Value = 1700
$ws.Range("M99").Value = -335.5714
$ws.Range("N99").Value = -4696

$ws.Range("H107").Value = 1574.1333
$ws.Range("I107").Value = 1599.9166
$ws.Range("K107").Value = 1599.9166
$ws.Range("M107").Value = 320.0834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 87000.39999999999
$ws.Range("I62").Value = 9665.666999999999
$ws.Range("J62").Value = 203002.5
$ws.Range("K62").Value = 9665.666999999999
$ws.Range("L62").Value = 203002.5
$ws.Range("M62").Value = -9041.666999999999
$ws.Range("N62").Value = -204250.5

$ws.Range("H65").Value = 87000.39999999999
$ws.Range("I65").Value = 9665.666999999999
$ws.Range("J65").Value = 203002.5
$ws.Range("K65").Value = 48328.335
$ws.Range("L65").Value = 1015012.5
$ws.Range("M65").Value = -45208.335
$ws.Range("N65").Value = -1021252.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1836.6765
$ws.Range("I11").Value = 1609.125
$ws.Range("K11").Value = 4827.375
$ws.Range("M11").Value = -4687.375

$ws.Range("H44").Value = 747.5
$ws.Range("I44").Value = 684.375
$ws.Range("J44").Value = 1000
$ws.Range("K44").Value = 2053.125
$ws.Range("L44").Value = 3000
$ws.Range("M44").Value = -1655.125
$ws.Range("N44").Value = -3796

$ws.Range("H68").Value = 1921.1515
$ws.Range("J68").Value = 2099.913
$ws.Range("L68").Value = 6299.739
$ws.Range("N68").Value = -7921.739

$ws.Range("H71").Value = 1921.1515
$ws.Range("J71").Value = 2099.913
$ws.Range("L71").Value = 18899.217
$ws.Range("N71").Value = -27011.217

$ws.Range("H131").Value = 1424.6
$ws.Range("I131").Value = 698
$ws.Range("J131").Value = 1439.4286
$ws.Range("K131").Value = 2094
$ws.Range("L131").Value = 4318.2858
$ws.Range("M131").Value = 2946
$ws.Range("N131").Value = -14398.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 195.54167
$ws.Range("I2").Value = 43
$ws.Range("J2").Value = 304.5
$ws.Range("K2").Value = 43
$ws.Range("L2").Value = 304.5
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = -530.5

$ws.Range("H122").Value = 58013.723
$ws.Range("I122").Value = 2602.8235
$ws.Range("K122").Value = 7808.470499999999
$ws.Range("M122").Value = -5358.470499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 992.5714
$ws.Range("I22").Value = 983
$ws.Range("K22").Value = 983
$ws.Range("M22").Value = -688

$ws.Range("H27").Value = 992.5714
$ws.Range("I27").Value = 983
$ws.Range("K27").Value = 983
$ws.Range("M27").Value = -876

$ws.Range("H42").Value = 1676666.6
$ws.Range("J42").Value = 1676666.6
$ws.Range("L42").Value = 1676666.6
$ws.Range("N42").Value = -1677792.6

$ws.Range("H43").Value = 3376116.5
$ws.Range("I43").Value = 10012
$ws.Range("J43").Value = 6069000
$ws.Range("K43").Value = 10012
$ws.Range("L43").Value = 6069000
$ws.Range("M43").Value = -9819
$ws.Range("N43").Value = -6069386

$ws.Range("H46").Value = 3249.75
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 3999.6667
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 3999.6667
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -4375.6667

$ws.Range("H49").Value = 1676666.6
$ws.Range("J49").Value = 1676666.6
$ws.Range("L49").Value = 1676666.6
$ws.Range("N49").Value = -1676960.6

$ws.Range("H122").Value = 7500
$ws.Range("J122").Value = 10000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -34900

$ws.Range("H125").Value = 95397.5
$ws.Range("J125").Value = 95397.5
$ws.Range("L125").Value = 95397.5
$ws.Range("N125").Value = -105237.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 15399.8
$ws.Range("J104").Value = 15399.8
$ws.Range("L104").Value = 15399.8
$ws.Range("N104").Value = -22387.8

$ws.Range("H122").Value = 3459.8
$ws.Range("I122").Value = 3100
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 9300
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -6850
$ws.Range("N122").Value = -16898.5

$ws.Range("H136").Value = 7771
$ws.Range("I136").Value = 7319.5
$ws.Range("J136").Value = 14995
$ws.Range("K136").Value = 21958.5
$ws.Range("L136").Value = 44985
$ws.Range("M136").Value = -19408.5
$ws.Range("N136").Value = -50085
